$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (D) column cells to keep their text representation
# (values like "538.52" would otherwise be auto-coerced to numbers).
$dCells = @("D2","D3","D5","D6","D9","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D28","D31","D32","D33","D36","D38","D39","D40","D41","D42","D43","D44","D46","D47","D48","D49","D50")
foreach ($ref in $dCells) { $ws.Range($ref).NumberFormat = "@" }

# Row 2
$ws.Range("D2").Value = "59.681.06"
$ws.Range("E2").Value = "  +0.94%  "

# Row 3
$ws.Range("D3").Value = "2.615.49"
$ws.Range("E3").Value = "  +0.99%  "

# Row 4
$ws.Range("E4").Value = "  -0.30%  "

# Row 5
$ws.Range("D5").Value = "538.52"
$ws.Range("E5").Value = "  +2.67%  "

# Row 6
$ws.Range("D6").Value = "142.32"
$ws.Range("E6").Value = "  +1.98%  "

# Row 7
$ws.Range("E7").Value = "  +0.26%  "

# Row 8
$ws.Range("E8").Value = "  +0.66%  "

# Row 9
$ws.Range("D9").Value = "6.57"
$ws.Range("E9").Value = "  +0.68%  "

# Row 10
$ws.Range("E10").Value = "  +1.16%  "

# Row 11
$ws.Range("E11").Value = "  +1.24%  "

# Row 12
$ws.Range("E12").Value = "  -1.19%  "

# Row 13
$ws.Range("D13").Value = "3.075.06"
$ws.Range("E13").Value = "  +0.87%  "

# Row 14
$ws.Range("D14").Value = "59.612.99"
$ws.Range("E14").Value = "  +0.97%  "

# Row 15
$ws.Range("D15").Value = "20.78"
$ws.Range("E15").Value = "  +1.38%  "

# Row 16
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000134"
$ws.Range("E16").Value = "  +0.64%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.587.94"
$ws.Range("E17").Value = "  -0.94%  "

# Row 18
$ws.Range("D18").Value = "340.97"
$ws.Range("E18").Value = "  -0.25%  "

# Row 19
$ws.Range("D19").Value = "4.36"
$ws.Range("E19").Value = "  +1.18%  "

# Row 20
$ws.Range("D20").Value = "10.13"
$ws.Range("E20").Value = "  +0.37%  "

# Row 21
$ws.Range("D21").Value = "6.36"
$ws.Range("E21").Value = "  -1.17%  "

# Row 22
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.11%  "

# Row 23
$ws.Range("D23").Value = "67.28"
$ws.Range("E23").Value = "  +0.88%  "

# Row 24
$ws.Range("D24").Value = "0.409"
$ws.Range("E24").Value = "  +0.85%  "

# Row 25
$ws.Range("E25").Value = "  -1.27%  "

# Row 26
$ws.Range("E26").Value = "  +0.15%  "

# Row 27
$ws.Range("E27").Value = "  +2.38%  "

# Row 28
$ws.Range("D28").Value = "0.0₃0749"
$ws.Range("E28").Value = "  +3.35%  "

# Row 29
$ws.Range("E29").Value = "  +0.04%  "

# Row 30
$ws.Range("E30").Value = "  +5.19%  "

# Row 31
$ws.Range("D31").Value = "5.84"
$ws.Range("E31").Value = "  -1.45%  "

# Row 32
$ws.Range("D32").Value = "18.83"
$ws.Range("E32").Value = "  +0.60%  "

# Row 33
$ws.Range("D33").Value = "150.73"
$ws.Range("E33").Value = "  +1.06%  "

# Row 34
$ws.Range("E34").Value = "  +0.58%  "

# Row 35
$ws.Range("E35").Value = "  +0.74%  "

# Row 36
$ws.Range("D36").Value = "0.836"
$ws.Range("E36").Value = "  +2.78%  "

# Row 37
$ws.Range("E37").Value = "  -0.99%  "

# Row 38
$ws.Range("D38").Value = "0.829"
$ws.Range("E38").Value = "  -0.13%  "

# Row 39
$ws.Range("D39").Value = "3.55"
$ws.Range("E39").Value = "  +0.77%  "

# Row 40
$ws.Range("D40").Value = "277.93"
$ws.Range("E40").Value = "  +2.38%  "

# Row 41
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  +0.26%  "

# Row 42
$ws.Range("D42").Value = "0.603"
$ws.Range("E42").Value = "  +1.09%  "

# Row 43
$ws.Range("D43").Value = "10.75"
$ws.Range("E43").Value = "  -0.22%  "

# Row 44
$ws.Range("D44").Value = "0.0950"
$ws.Range("E44").Value = "  -0.08%  "

# Row 45
$ws.Range("E45").Value = "  +1.89%  "

# Row 46
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "1.948.60"
$ws.Range("E46").Value = "  -1.22%  "

# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0224"
$ws.Range("E47").Value = "  +0.45%  "

# Row 48
$ws.Range("D48").Value = "18.52"
$ws.Range("E48").Value = "  +1.93%  "

# Row 49
$ws.Range("D49").Value = "4.52"
$ws.Range("E49").Value = "  +1.17%  "

# Row 50
$ws.Range("D50").Value = "110.72"
$ws.Range("E50").Value = "  -3.40%  "

# Row 51
$ws.Range("E51").Value = "  +0.90%  "

# Restore the Normal style on those cells so no stray text-format style sticks
foreach ($ref in $dCells) { $ws.Range($ref).Style = "Normal" }